$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# Sheet1: new "multi turn" block (rows 11-17), mirrors the existing
# "single turn" block (rows 3-9) in layout and styling.
# ------------------------------------------------------------------

# Row 11: section title "multi turn" (bold label style, like B5/B6/B8/B9)
$ws1.Range("B11").Value = "multi turn"
$ws1.Range("B5").Copy()
$ws1.Range("B11").PasteSpecial(-4122)

# Row 12: header row (gpt3 wins / gpt2 wins / tie / d / sum), mirrors row 4
$ws1.Range("B12").Value = ""
$ws1.Range("B4").Copy()
$ws1.Range("B12").PasteSpecial(-4122)

$ws1.Range("C12").Value = "gpt3 wins"
$ws1.Range("C4").Copy()
$ws1.Range("C12").PasteSpecial(-4122)

$ws1.Range("D12").Value = "gpt2 wins"
$ws1.Range("D4").Copy()
$ws1.Range("D12").PasteSpecial(-4122)

$ws1.Range("E12").Value = "tie"
$ws1.Range("E4").Copy()
$ws1.Range("E12").PasteSpecial(-4122)

$ws1.Range("F12").Value = "sum"
$ws1.Range("F4").Copy()
$ws1.Range("F12").PasteSpecial(-4122)

$ws1.Range("G12").Value = "d"
$ws1.Range("G4").Copy()
$ws1.Range("G12").PasteSpecial(-4122)

# Row 13: naturalness data row, mirrors row 5 (but F is a SUM formula)
$ws1.Range("B13").Value = "naturalness"
$ws1.Range("B5").Copy()
$ws1.Range("B13").PasteSpecial(-4122)

$ws1.Range("C13").Value = 40
$ws1.Range("C5").Copy()
$ws1.Range("C13").PasteSpecial(-4122)

$ws1.Range("D13").Value = 16
$ws1.Range("D5").Copy()
$ws1.Range("D13").PasteSpecial(-4122)

$ws1.Range("E13").Value = 33
$ws1.Range("E5").Copy()
$ws1.Range("E13").PasteSpecial(-4122)

$ws1.Range("F13").Formula = "=SUM(C13:E13)"
$ws1.Range("F5").Copy()
$ws1.Range("F13").PasteSpecial(-4122)

$ws1.Range("G13").Formula = "=C13-D13"
$ws1.Range("H13").FormulaArray = "=[1]!TRINOM_TEST(G13/2,E13/2,F13/2,1)"

# Row 14: usefulness data row, mirrors row 6 (shared formula for G)
$ws1.Range("B14").Value = "usefulness"
$ws1.Range("B6").Copy()
$ws1.Range("B14").PasteSpecial(-4122)

$ws1.Range("C14").Value = 55
$ws1.Range("C6").Copy()
$ws1.Range("C14").PasteSpecial(-4122)

$ws1.Range("D14").Value = 11
$ws1.Range("D6").Copy()
$ws1.Range("D14").PasteSpecial(-4122)

$ws1.Range("E14").Value = 23
$ws1.Range("E6").Copy()
$ws1.Range("E14").PasteSpecial(-4122)

$ws1.Range("F14").Formula = "=SUM(C14:E14)"
$ws1.Range("F6").Copy()
$ws1.Range("F14").PasteSpecial(-4122)

$ws1.Range("G14").Formula = "=C14-D14"
$ws1.Range("H14").FormulaArray = "=[1]!TRINOM_TEST(G14/2,E14/2,F14/2,1)"

# Row 15: second header row (gpt3 wins / human wins / tie), mirrors row 7
$ws1.Range("B15").Value = ""
$ws1.Range("B7").Copy()
$ws1.Range("B15").PasteSpecial(-4122)

$ws1.Range("C15").Value = "gpt3 wins"
$ws1.Range("C7").Copy()
$ws1.Range("C15").PasteSpecial(-4122)

$ws1.Range("D15").Value = "human wins"
$ws1.Range("D7").Copy()
$ws1.Range("D15").PasteSpecial(-4122)

$ws1.Range("E15").Value = "tie"
$ws1.Range("E7").Copy()
$ws1.Range("E15").PasteSpecial(-4122)

$ws1.Range("F15").Value = ""
$ws1.Range("F7").Copy()
$ws1.Range("F15").PasteSpecial(-4122)

# Row 16: naturalness data row, mirrors row 8
$ws1.Range("B16").Value = "naturalness"
$ws1.Range("B8").Copy()
$ws1.Range("B16").PasteSpecial(-4122)

$ws1.Range("C16").Value = 16
$ws1.Range("C8").Copy()
$ws1.Range("C16").PasteSpecial(-4122)

$ws1.Range("D16").Value = 18
$ws1.Range("D8").Copy()
$ws1.Range("D16").PasteSpecial(-4122)

$ws1.Range("E16").Value = 31
$ws1.Range("E8").Copy()
$ws1.Range("E16").PasteSpecial(-4122)

$ws1.Range("F16").Formula = "=SUM(C16:E16)"
$ws1.Range("F8").Copy()
$ws1.Range("F16").PasteSpecial(-4122)

$ws1.Range("G16:G17").Formula = "=C16-D16"
$ws1.Range("H16").FormulaArray = "=[1]!TRINOM_TEST(G16/2,E16/2,F16/2,1)"
$ws1.Range("I16").Value = "swap 16 and 18 for p-value"

# Row 17: usefulness data row, mirrors row 9
$ws1.Range("B17").Value = "usefulness"
$ws1.Range("B9").Copy()
$ws1.Range("B17").PasteSpecial(-4122)

$ws1.Range("C17").Value = 17
$ws1.Range("C9").Copy()
$ws1.Range("C17").PasteSpecial(-4122)

$ws1.Range("D17").Value = 10
$ws1.Range("D9").Copy()
$ws1.Range("D17").PasteSpecial(-4122)

$ws1.Range("E17").Value = 38
$ws1.Range("E9").Copy()
$ws1.Range("E17").PasteSpecial(-4122)

$ws1.Range("F17").Formula = "=SUM(C17:E17)"
$ws1.Range("F9").Copy()
$ws1.Range("F17").PasteSpecial(-4122)

$ws1.Range("H17").FormulaArray = "=[1]!TRINOM_TEST(G17/2,E17/2,F17/2,1)"

# ------------------------------------------------------------------
# Add Sheet2 (Assignment / Score / Perc table), placed after Sheet1
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Assignment"
$ws2.Range("B1").Value = "Score"
$ws2.Range("C1").Value = "Perc"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 63
$ws2.Range("C2").Value = 0.15
$ws2.Range("D2").Formula = "=B2*C2"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 99.5
$ws2.Range("C3").Value = 0.2
$ws2.Range("D3:D5").Formula = "=B3*C3"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 78
$ws2.Range("C4").Value = 0.25

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = 24
$ws2.Range("C5").Value = 0.4

$ws2.Range("D6").Formula = "=SUM(D2:D5)"

$ws2.Range("K9").Select()

# ------------------------------------------------------------------
# Back on Sheet1: "Oleg:" note plus a couple more trinomial test rows
# ------------------------------------------------------------------
$ws1.Range("B22").Value = "Oleg:"

$ws1.Range("C23").Value = 60
$ws1.Range("D23").Value = 16
$ws1.Range("E23").Value = 24
$ws1.Range("F23").Formula = "=SUM(C23:E23)"
$ws1.Range("G23").Formula = "=C23-D23"
$ws1.Range("H23").FormulaArray = "=[1]!TRINOM_TEST(G23/2,E23/2,F23/2,1)"

$ws1.Range("C24").Value = 29
$ws1.Range("D24").Value = 28
$ws1.Range("E24").Value = 43
$ws1.Range("F24").Formula = "=SUM(C24:E24)"
$ws1.Range("G24").Formula = "=C24-D24"
$ws1.Range("H24").FormulaArray = "=[1]!TRINOM_TEST(G24/2,E24/2,F24/2,1)"

$ws1.Range("F24").Select()
